$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B for the new "Priority" field.
# This shifts the old "Groups" column (B) to C and the old "Execution" column (C) to D.
$ws.Columns.Item(2).Insert()

# Headers
$ws.Range("B1").Value = "Priority"
$ws.Range("C1").Value = "Group"

# Priority values for each test case row (TestCaseNumber 101-111).
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 4
$ws.Range("B12").Value = 4

# The old combined "Groups=Smoke,Regression" value shifted into D2; replace it with "All".
$ws.Range("D2").Value = "All"

# Match formatting: TestCaseNumber & Priority columns are left-aligned, stored as Text.
$ws.Range("A1:B12").NumberFormat = "@"
$ws.Range("A1:B12").HorizontalAlignment = -4131

# Match column widths: Priority column uses the same width as TestCaseNumber.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

$ws.Range("D4").Select()
